# Applies the "Data Munging & Plot Work" edits to the Chart To-Do List workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 ("Total Carbon Offset by State"): Completed? -> "Yes" (was "In Progress")
$ws.Range("E5").Value = "Yes"
$ws.Range("E5").Style = "Good"

# Column I ("Completed?3" - Adjusted Carbon Offset related rows) for rows 6-14:
# was "No" -> now "In Progress"
$rng = $ws.Range("I6:I14")
$rng.Value = "In Progress"
$rng.Style = "Neutral"

# Update the active selection to reflect where the user ended up (I18)
$ws.Range("I18").Select()
